$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.214453472130288
